$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "matcher" column (old C),
# pushing it to column E and opening up C:D for the new "elementType"/
# "Locator" columns.
$ws.Columns("C:D").Insert()

# Fill the new "elementType" column data first (matches the order the
# shared-string table ends up in), then the header row.
$ws.Range("C2").Value = "textbox"
$ws.Range("C3").Value = "textbox"
$ws.Range("C4").Value = "button"
$ws.Range("C1").Value = "elementType"
$ws.Range("D1").Value = "Locator"

# The insert leaves a blank styled placeholder at D2 (row 2 carries an
# explicit row style); strip it back down to a plain empty cell.
$ws.Range("D2").Style = "Normal"

# Match the column widths used for the new columns.
$ws.Columns("C").ColumnWidth = 12.25
$ws.Columns("D").ColumnWidth = 6.6

# Update the selected cell to follow the last populated cell (E4).
$ws.Range("E4").Select() | Out-Null
